# Mise à jour de l'application
# Fill down the attendance-summary formulas (columns B..J) from row 26
# into rows 27-29 (COUNTA/COUNTIF totals for the three remaining players),
# and update the frozen-pane view / active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the actual COUNTA/COUNTIF summary formulas for each new row first...
for ($r = 27; $r -le 29; $r++) {
    $ws.Cells.Item($r, 2).Formula = "=COUNTA(K$r`:VQ$r)"
    $ws.Cells.Item($r, 3).Formula = "=COUNTIF(K$r`:VQ$r,`"P`")"
    $ws.Cells.Item($r, 4).Formula = "=COUNTIF(K$r`:VQ$r,`"REP`")"
    $ws.Cells.Item($r, 5).Formula = "=COUNTIF(K$r`:VR$r,`"A`")"
    $ws.Cells.Item($r, 6).Formula = "=COUNTIF(K$r`:VS$r,`"B`")"
    $ws.Cells.Item($r, 7).Formula = "=COUNTIF(K$r`:VT$r,`"M`")"
    $ws.Cells.Item($r, 8).Formula = "=COUNTIF(K$r`:VU$r,`"R`")"
    $ws.Cells.Item($r, 9).Formula = "=COUNTIF(K$r`:VV$r,`"RH`")"
    $ws.Cells.Item($r, 10).Formula = "=COUNTIF(K$r`:VV$r,`"S`")"
}

# ...then copy the formatting from the last filled-in summary row (26) down
# onto the three new rows so the cell styles match the rest of the column.
$ws.Range("B26:J26").Copy()
$ws.Range("B27:J29").PasteSpecial(-4122)

# Update the view: scroll the frozen pane back to column B and move the
# active selection to L30.
$ws.Range("L30").Select()
